$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price cells so the exact
# source formatting (trailing zeros, etc.) is preserved as text,
# matching the workbook's original inlineStr string cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the crypto price refresh.
$ws.Range("D2").Value = '26.666.33'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '1.794.09'
$ws.Range("E3").Value = '  -1.34%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '308.97'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '0.4441'
$ws.Range("E7").Value = '  +5.42%  '
$ws.Range("D8").Value = '0.3676'
$ws.Range("E8").Value = '  +0.38%  '
$ws.Range("D9").Value = '0.07333'
$ws.Range("E9").Value = '  +1.85%  '
$ws.Range("D10").Value = '0.8561'
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("D11").Value = '20.61'
$ws.Range("E11").Value = '  -0.95%  '
$ws.Range("D12").Value = '1.799.44'
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("D13").Value = '6.613'
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("D14").Value = '91.95'
$ws.Range("E14").Value = '  +3.55%  '
$ws.Range("D15").Value = '0.07063'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").Value = '5.259'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '0.000008671'
$ws.Range("E18").Value = '  -1.34%  '
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '14.77'
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").Value = '26.688.40'
$ws.Range("E21").Value = '  -1.50%  '
$ws.Range("D22").Value = '5.150'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").Value = '1.981'
$ws.Range("E24").Value = '  +0.33%  '
$ws.Range("D25").Value = '151.84'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.182'
$ws.Range("E26").Value = '  -2.37%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '18.36'
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").Value = '5.166'
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("D29").Value = '117.20'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("D30").Value = '0.08760'
$ws.Range("E30").Value = '  -0.25%  '
$ws.Range("D31").Value = '0.7369'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").Value = '1.154'
$ws.Range("E32").Value = '  -1.77%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.433'
$ws.Range("E33").Value = '  +0.62%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.901'
$ws.Range("E34").Value = '  -2.88%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("D36").Value = '1.084'
$ws.Range("E36").Value = '  -1.33%  '
$ws.Range("D37").Value = '0.01952'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").Value = '0.05174'
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("D39").Value = '0.5243'
$ws.Range("E39").Value = '  +4.50%  '
$ws.Range("D40").Value = '2.823'
$ws.Range("E40").Value = '  -1.85%  '
$ws.Range("D41").Value = '6.961'
$ws.Range("E41").Value = '  -4.08%  '
$ws.Range("D42").Value = '0.1679'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("D43").Value = '0.5039'
$ws.Range("E43").Value = '  +6.48%  '
$ws.Range("D44").Value = '8.418'
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").Value = '1.965'
$ws.Range("E45").Value = '  +4.70%  '
$ws.Range("D46").Value = '10.45'
$ws.Range("E46").Value = '  -0.84%  '
$ws.Range("D47").Value = '104.94'
$ws.Range("E47").Value = '  -1.03%  '
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").Value = '1.663'
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("D50").Value = '0.06288'
$ws.Range("D51").Value = '0.9135'
$ws.Range("E51").Value = '  +1.38%  '
